# coutOptimalParInstance.xlsx update:
#  - add a 7th instance row (A8/B8) with its computed cost
#  - refresh the computed costs for instance5 (B6) and instance6 (B7)
#  - move the active selection to B11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated costs for existing instances
$ws.Range("B6").Value = 12.504761904761899
$ws.Range("B7").Value = 1.2523809523809499

# New instance7 row
$ws.Range("A8").Value = "instance7"
$ws.Range("B8").Value = 49.884303350970001
# instance7's cost cell has no surrounding border, same look as the other
# "no border" value cells in this column
$ws.Range("B8").Borders.LineStyle = -4142

# The border-less formatting that used to single out B3 is no longer needed
$ws.Range("B3").Style = "Normal"

# Move / restore the selection
$null = $ws.Range("B11").Select()
